$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Remove the second worker's record (Camila Marcela Primera Guerrero) entirely
$ws.Rows("17:17").Delete()

# Update the "Periodo Mora" value for the remaining worker (Liliana) from 2507 to 2508
$ws.Range("E16").Value = "2508"

# Update totals to reflect the removed worker
$ws.Range("E11").Value = 56940
$ws.Range("C13").Value = 1
